$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.915.14"
$ws.Range("E2").Value = "  +0.89%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.924.27"
$ws.Range("E3").Value = "  -0.40%  "

# Row 4
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.52"
$ws.Range("E5").Value = "  +6.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.56"
$ws.Range("E6").Value = "  -1.69%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  -1.80%  "

# Row 8
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.726"
$ws.Range("E9").Value = "  -0.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  -2.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000342"
$ws.Range("E11").Value = "  -2.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.20"
$ws.Range("E12").Value = "  -2.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.47"
$ws.Range("E13").Value = "  +0.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.522.67"
$ws.Range("E14").Value = "  -1.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.908.16"
$ws.Range("E15").Value = "  -0.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.31"
$ws.Range("E16").Value = "  +7.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.16"
$ws.Range("E17").Value = "  -0.55%  "

# Row 18
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.21"
$ws.Range("E18").Value = "  +2.46%  "

# Row 19
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.134"
$ws.Range("E19").Value = "  -1.86%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.710.81"
$ws.Range("E20").Value = "  +0.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "423.94"
$ws.Range("E21").Value = "  -3.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("E22").Value = "  +1.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.12"
$ws.Range("E23").Value = "  -3.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.82"
$ws.Range("E24").Value = "  -1.36%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.04"
$ws.Range("E25").Value = "  +6.40%  "

# Row 26
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.87"
$ws.Range("E26").Value = "  -1.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.60"
$ws.Range("E27").Value = "  -4.63%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.13"
$ws.Range("E28").Value = "  -2.91%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.45"
$ws.Range("E29").Value = "  +0.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "685.64"
$ws.Range("E30").Value = "  -2.31%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.06"
$ws.Range("E31").Value = "  +16.40%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.126"
$ws.Range("E32").Value = "  -3.36%  "

# Row 33
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.87"
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.91"
$ws.Range("E34").Value = "  +7.91%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.441"
$ws.Range("E35").Value = "  -4.67%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0879"
$ws.Range("E36").Value = "  -1.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.86"
$ws.Range("E37").Value = "  -2.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.41"
$ws.Range("E38").Value = "  +11.32%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.148"
$ws.Range("E39").Value = "  -1.60%  "

# Row 40
$ws.Range("E40").Value = "  +0.23%  "

# Row 41
$ws.Range("E41").Value = "  -0.45%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0480"
$ws.Range("E42").Value = "  -1.71%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.16"
$ws.Range("E43").Value = "  +4.59%  "

# Row 44
$ws.Range("E44").Value = "  -2.29%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.54"
$ws.Range("E45").Value = "  +4.42%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.143"
$ws.Range("E46").Value = "  -0.47%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.12"
$ws.Range("E47").Value = "  +4.14%  "

# Row 48
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000277"
$ws.Range("E48").Value = "  +16.93%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.35"
$ws.Range("E49").Value = "  -1.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.91"
$ws.Range("E50").Value = "  +2.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.98"
$ws.Range("E51").Value = "  -0.43%  "
